$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 3371.2856  # H18: 3099.8333 -> 3371.2856
$ws.Cells.Item(18, 9).Value = 1619.8  # I18: 774.75 -> 1619.8
$ws.Cells.Item(18, 11).Value = 1619.8  # K18: 774.75 -> 1619.8
$ws.Cells.Item(18, 13).Value = -1335.8  # M18: -490.75 -> -1335.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 6195.8887  # H69: 6293.8335 -> 6195.8887
$ws.Cells.Item(69, 9).Value = 4006.5  # I69: 4013 -> 4006.5
$ws.Cells.Item(69, 10).Value = 6821.4287  # J69: 6750 -> 6821.4287
$ws.Cells.Item(69, 11).Value = 12019.5  # K69: 12039 -> 12019.5
$ws.Cells.Item(69, 12).Value = 20464.2861  # L69: 20250 -> 20464.2861
$ws.Cells.Item(69, 13).Value = -11145.5  # M69: -11165 -> -11145.5
$ws.Cells.Item(69, 14).Value = -22212.2861  # N69: -21998 -> -22212.2861

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 6195.8887  # H72: 6293.8335 -> 6195.8887
$ws.Cells.Item(72, 9).Value = 4006.5  # I72: 4013 -> 4006.5
$ws.Cells.Item(72, 10).Value = 6821.4287  # J72: 6750 -> 6821.4287
$ws.Cells.Item(72, 11).Value = 36058.5  # K72: 36117 -> 36058.5
$ws.Cells.Item(72, 12).Value = 61392.85830000001  # L72: 60750 -> 61392.85830000001
$ws.Cells.Item(72, 13).Value = -31690.5  # M72: -31749 -> -31690.5
$ws.Cells.Item(72, 14).Value = -70128.85830000001  # N72: -69486 -> -70128.85830000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 4038.0557  # H98: 3938.5945 -> 4038.0557
$ws.Cells.Item(98, 9).Value = 3572.375  # I98: 3474.9697 -> 3572.375
$ws.Cells.Item(98, 11).Value = 3572.375  # K98: 3474.9697 -> 3572.375
$ws.Cells.Item(98, 13).Value = -2074.375  # M98: -1976.9697 -> -2074.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 2285.3103  # H106: 2132.96 -> 2285.3103
$ws.Cells.Item(106, 10).Value = 2399.8333  # J106: 724.5 -> 2399.8333
$ws.Cells.Item(106, 12).Value = 2399.8333  # L106: 724.5 -> 2399.8333
$ws.Cells.Item(106, 14).Value = -3661.8333  # N106: -1986.5 -> -3661.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 4038.0557  # H122: 3938.5945 -> 4038.0557
$ws.Cells.Item(122, 9).Value = 3572.375  # I122: 3474.9697 -> 3572.375
$ws.Cells.Item(122, 11).Value = 10717.125  # K122: 10424.9091 -> 10717.125
$ws.Cells.Item(122, 13).Value = -8267.125  # M122: -7974.909100000001 -> -8267.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1845.439  # H132: 1891.325 -> 1845.439
$ws.Cells.Item(132, 9).Value = 1516.8823  # I132: 1592.6562 -> 1516.8823
$ws.Cells.Item(132, 10).Value = 3441.2856  # J132: 3086 -> 3441.2856
$ws.Cells.Item(132, 11).Value = 4550.6469  # K132: 4777.9686 -> 4550.6469
$ws.Cells.Item(132, 12).Value = 10323.8568  # L132: 9258 -> 10323.8568
$ws.Cells.Item(132, 13).Value = -2020.6469  # M132: -2247.9686 -> -2020.6469
$ws.Cells.Item(132, 14).Value = -15383.8568  # N132: -14318 -> -15383.8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2195.3333  # H137: 2407.3333 -> 2195.3333
$ws.Cells.Item(137, 9).Value = 1573.8182  # I137: 1726.8889 -> 1573.8182
$ws.Cells.Item(137, 10).Value = 2879  # J137: 3087.7778 -> 2879
$ws.Cells.Item(137, 11).Value = 4721.4546  # K137: 5180.6667 -> 4721.4546
$ws.Cells.Item(137, 12).Value = 8637  # L137: 9263.3334 -> 8637
$ws.Cells.Item(137, 13).Value = -2171.4546  # M137: -2630.6667 -> -2171.4546
$ws.Cells.Item(137, 14).Value = -13737  # N137: -14363.3334 -> -13737

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2884.88  # H138: 2959.61 -> 2884.88
$ws.Cells.Item(138, 9).Value = 1893.9642  # I138: 2110.7 -> 1893.9642
$ws.Cells.Item(138, 10).Value = 3270.236  # J138: 3323.4285 -> 3270.236
$ws.Cells.Item(138, 11).Value = 5681.892599999999  # K138: 6332.099999999999 -> 5681.892599999999
$ws.Cells.Item(138, 12).Value = 9810.707999999999  # L138: 9970.2855 -> 9810.707999999999
$ws.Cells.Item(138, 13).Value = -541.8925999999992  # M138: -1192.099999999999 -> -541.8925999999992
$ws.Cells.Item(138, 14).Value = -20090.708  # N138: -20250.2855 -> -20090.708

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 2501  # H14: 2502 -> 2501
$ws.Cells.Item(14, 10).Value = 2501.6667  # J14: 2503.5 -> 2501.6667
$ws.Cells.Item(14, 12).Value = 2501.6667  # L14: 2503.5 -> 2501.6667
$ws.Cells.Item(14, 14).Value = -2851.6667  # N14: -2853.5 -> -2851.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 999.3333  # H17: 799.6 -> 999.3333
$ws.Cells.Item(17, 9).Value = 0  # I17: 500 -> 0
$ws.Cells.Item(17, 11).Value = 0  # K17: 500 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # M17 remove (was -327)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 5558016.5  # H122: 5850464 -> 5558016.5
$ws.Cells.Item(122, 9).Value = 6538313.5  # I122: 6946863.5 -> 6538313.5
$ws.Cells.Item(122, 11).Value = 19614940.5  # K122: 20840590.5 -> 19614940.5
$ws.Cells.Item(122, 13).Value = -19612490.5  # M122: -20838140.5 -> -19612490.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 5500.8887  # H132: 6356.9565 -> 5500.8887
$ws.Cells.Item(132, 9).Value = 5569.88  # I132: 6520.619 -> 5569.88
$ws.Cells.Item(132, 11).Value = 16709.64  # K132: 19561.857 -> 16709.64
$ws.Cells.Item(132, 13).Value = -14179.64  # M132: -17031.857 -> -14179.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 5666.6665  # H16: 0 -> 5666.6665
$ws.Cells.Item(16, 9).Value = 5000  # I16: 0 -> 5000
$ws.Cells.Item(16, 10).Value = 7000  # J16: 0 -> 7000
$ws.Cells.Item(16, 11).Value = 5000  # K16: 0 -> 5000
$ws.Cells.Item(16, 12).Value = 7000  # L16: 0 -> 7000
$ws.Cells.Item(16, 13).Value = -4830  # M16: None -> -4830
$ws.Cells.Item(16, 14).Value = -7340  # N16: None -> -7340

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2305.4285  # H94: 2306.2144 -> 2305.4285
$ws.Cells.Item(94, 9).Value = 2305.4285  # I94: 2375.923 -> 2305.4285
$ws.Cells.Item(94, 10).Value = 0  # J94: 1400 -> 0
$ws.Cells.Item(94, 11).Value = 2305.4285  # K94: 2375.923 -> 2305.4285
$ws.Cells.Item(94, 12).Value = 0  # L94: 1400 -> 0
$ws.Cells.Item(94, 13).Value = -1854.4285  # M94: -1924.923 -> -1854.4285
$ws.Cells.Item(94, 14).ClearContents()  # N94 remove (was -2302)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3911.6785  # H99: 3956.2307 -> 3911.6785
$ws.Cells.Item(99, 9).Value = 3651.4707  # I99: 3692.25 -> 3651.4707
$ws.Cells.Item(99, 10).Value = 4313.8184  # J99: 4378.6 -> 4313.8184
$ws.Cells.Item(99, 11).Value = 3651.4707  # K99: 3692.25 -> 3651.4707
$ws.Cells.Item(99, 12).Value = 4313.8184  # L99: 4378.6 -> 4313.8184
$ws.Cells.Item(99, 13).Value = -2153.4707  # M99: -2194.25 -> -2153.4707
$ws.Cells.Item(99, 14).Value = -7309.8184  # N99: -7374.6 -> -7309.8184

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4384.718  # H107: 4473.921 -> 4384.718
$ws.Cells.Item(107, 9).Value = 3240.1667  # I107: 3317.5862 -> 3240.1667
$ws.Cells.Item(107, 11).Value = 3240.1667  # K107: 3317.5862 -> 3240.1667
$ws.Cells.Item(107, 13).Value = -1320.1667  # M107: -1397.5862 -> -1320.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1246.1666  # H16: 1242.8462 -> 1246.1666
$ws.Cells.Item(16, 9).Value = 1048.5  # I16: 1009.5 -> 1048.5
$ws.Cells.Item(16, 10).Value = 1443.8334  # J16: 1616.2 -> 1443.8334
$ws.Cells.Item(16, 11).Value = 1048.5  # K16: 1009.5 -> 1048.5
$ws.Cells.Item(16, 12).Value = 1443.8334  # L16: 1616.2 -> 1443.8334
$ws.Cells.Item(16, 13).Value = -761.5  # M16: -722.5 -> -761.5
$ws.Cells.Item(16, 14).Value = -2017.8334  # N16: -2190.2 -> -2017.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2956.7874  # H31: 2736.5283 -> 2956.7874
$ws.Cells.Item(31, 9).Value = 1967.05  # I31: 1842.3695 -> 1967.05
$ws.Cells.Item(31, 11).Value = 1967.05  # K31: 1842.3695 -> 1967.05
$ws.Cells.Item(31, 13).Value = -1672.05  # M31: -1547.3695 -> -1672.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2956.7874  # H34: 2736.5283 -> 2956.7874
$ws.Cells.Item(34, 9).Value = 1967.05  # I34: 1842.3695 -> 1967.05
$ws.Cells.Item(34, 11).Value = 1967.05  # K34: 1842.3695 -> 1967.05
$ws.Cells.Item(34, 13).Value = -1765.05  # M34: -1640.3695 -> -1765.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 8190.1665  # H62: 8498.362999999999 -> 8190.1665
$ws.Cells.Item(62, 9).Value = 7785.625  # I62: 8212.143 -> 7785.625
$ws.Cells.Item(62, 11).Value = 7785.625  # K62: 8212.143 -> 7785.625
$ws.Cells.Item(62, 13).Value = -7161.625  # M62: -7588.143 -> -7161.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 8190.1665  # H65: 8498.362999999999 -> 8190.1665
$ws.Cells.Item(65, 9).Value = 7785.625  # I65: 8212.143 -> 7785.625
$ws.Cells.Item(65, 11).Value = 38928.125  # K65: 41060.715 -> 38928.125
$ws.Cells.Item(65, 13).Value = -35808.125  # M65: -37940.715 -> -35808.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 55555.445  # H68: 52777.668 -> 55555.445
$ws.Cells.Item(68, 10).Value = 57499.832  # J68: 53333.168 -> 57499.832
$ws.Cells.Item(68, 12).Value = 57499.832  # L68: 53333.168 -> 57499.832
$ws.Cells.Item(68, 14).Value = -58997.832  # N68: -54831.168 -> -58997.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 55555.445  # H71: 52777.668 -> 55555.445
$ws.Cells.Item(71, 10).Value = 57499.832  # J71: 53333.168 -> 57499.832
$ws.Cells.Item(71, 12).Value = 172499.496  # L71: 159999.504 -> 172499.496
$ws.Cells.Item(71, 14).Value = -179987.496  # N71: -167487.504 -> -179987.496

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 2581.08  # H105: 2581.48 -> 2581.08
$ws.Cells.Item(105, 9).Value = 3868.25  # I105: 3625.3076 -> 3868.25
$ws.Cells.Item(105, 10).Value = 1392.9231  # J105: 1450.6666 -> 1392.9231
$ws.Cells.Item(105, 11).Value = 3868.25  # K105: 3625.3076 -> 3868.25
$ws.Cells.Item(105, 12).Value = 1392.9231  # L105: 1450.6666 -> 1392.9231
$ws.Cells.Item(105, 13).Value = -2121.25  # M105: -1878.3076 -> -2121.25
$ws.Cells.Item(105, 14).Value = -4886.9231  # N105: -4944.6666 -> -4886.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1246.1666  # H113: 1242.8462 -> 1246.1666
$ws.Cells.Item(113, 9).Value = 1048.5  # I113: 1009.5 -> 1048.5
$ws.Cells.Item(113, 10).Value = 1443.8334  # J113: 1616.2 -> 1443.8334
$ws.Cells.Item(113, 11).Value = 1048.5  # K113: 1009.5 -> 1048.5
$ws.Cells.Item(113, 12).Value = 1443.8334  # L113: 1616.2 -> 1443.8334
$ws.Cells.Item(113, 13).Value = 1121.5  # M113: 1160.5 -> 1121.5
$ws.Cells.Item(113, 14).Value = -5783.8334  # N113: -5956.2 -> -5783.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 4337.5835  # H122: 4279.1665 -> 4337.5835
$ws.Cells.Item(122, 9).Value = 4504.636  # I122: 4279.1665 -> 4504.636
$ws.Cells.Item(122, 10).Value = 2500  # J122: 0 -> 2500
$ws.Cells.Item(122, 11).Value = 13513.908  # K122: 12837.4995 -> 13513.908
$ws.Cells.Item(122, 12).Value = 7500  # L122: 0 -> 7500
$ws.Cells.Item(122, 13).Value = -11063.908  # M122: -10387.4995 -> -11063.908
$ws.Cells.Item(122, 14).Value = -12400  # N122: None -> -12400

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2534  # H5: 2571.2222 -> 2534
$ws.Cells.Item(5, 9).Value = 2295.2  # I5: 2319.25 -> 2295.2
$ws.Cells.Item(5, 11).Value = 6885.599999999999  # K5: 6957.75 -> 6885.599999999999
$ws.Cells.Item(5, 13).Value = -6773.599999999999  # M5: -6845.75 -> -6773.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 2715.8462  # H23: 3221.4666 -> 2715.8462
$ws.Cells.Item(23, 10).Value = 1333.1  # J23: 2195.5833 -> 1333.1
$ws.Cells.Item(23, 12).Value = 3999.3  # L23: 6586.749899999999 -> 3999.3
$ws.Cells.Item(23, 14).Value = -4469.299999999999  # N23: -7056.749899999999 -> -4469.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2669.2083  # H132: 2720.087 -> 2669.2083
$ws.Cells.Item(132, 10).Value = 2675.5908  # J132: 2731.6191 -> 2675.5908
$ws.Cells.Item(132, 12).Value = 24080.3172  # L132: 24584.5719 -> 24080.3172
$ws.Cells.Item(132, 14).Value = -29140.3172  # N132: -29644.5719 -> -29140.3172

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 2534  # H135: 2571.2222 -> 2534
$ws.Cells.Item(135, 9).Value = 2295.2  # I135: 2319.25 -> 2295.2
$ws.Cells.Item(135, 11).Value = 20656.8  # K135: 20873.25 -> 20656.8
$ws.Cells.Item(135, 13).Value = -18121.8  # M135: -18338.25 -> -18121.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 225  # H41: 1305.5 -> 225
$ws.Cells.Item(41, 9).Value = 200  # I41: 1305.5 -> 200
$ws.Cells.Item(41, 10).Value = 250  # J41: 0 -> 250
$ws.Cells.Item(41, 11).Value = 200  # K41: 1305.5 -> 200
$ws.Cells.Item(41, 12).Value = 250  # L41: 0 -> 250
$ws.Cells.Item(41, 13).Value = 155  # M41: -950.5 -> 155
$ws.Cells.Item(41, 14).Value = -960  # N41: None -> -960

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 4253.0527  # H46: 3641.3333 -> 4253.0527
$ws.Cells.Item(46, 9).Value = 4253.0527  # I46: 3641.3333 -> 4253.0527
$ws.Cells.Item(46, 11).Value = 4253.0527  # K46: 3641.3333 -> 4253.0527
$ws.Cells.Item(46, 13).Value = -4097.0527  # M46: -3485.3333 -> -4097.0527

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 37000  # H49: 0 -> 37000
$ws.Cells.Item(49, 10).Value = 37000  # J49: 0 -> 37000
$ws.Cells.Item(49, 12).Value = 37000  # L49: 0 -> 37000
$ws.Cells.Item(49, 14).Value = -37368  # N49: None -> -37368

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 25645734  # H70: 27782630 -> 25645734
$ws.Cells.Item(70, 10).Value = 5166.6665  # J70: 6250 -> 5166.6665
$ws.Cells.Item(70, 12).Value = 5166.6665  # L70: 6250 -> 5166.6665
$ws.Cells.Item(70, 14).Value = -5706.6665  # N70: -6790 -> -5706.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 25645734  # H73: 27782630 -> 25645734
$ws.Cells.Item(73, 10).Value = 5166.6665  # J73: 6250 -> 5166.6665
$ws.Cells.Item(73, 12).Value = 5166.6665  # L73: 6250 -> 5166.6665
$ws.Cells.Item(73, 14).Value = -7038.6665  # N73: -8122 -> -7038.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1073.5294  # H97: 1056.2727 -> 1073.5294
$ws.Cells.Item(97, 9).Value = 1046.875  # I97: 1056.2727 -> 1046.875
$ws.Cells.Item(97, 10).Value = 1500  # J97: 0 -> 1500
$ws.Cells.Item(97, 11).Value = 1046.875  # K97: 1056.2727 -> 1046.875
$ws.Cells.Item(97, 12).Value = 1500  # L97: 0 -> 1500
$ws.Cells.Item(97, 13).Value = -550.875  # M97: -560.2727 -> -550.875
$ws.Cells.Item(97, 14).Value = -2492  # N97: None -> -2492

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6174.9473  # H122: 6066.2 -> 6174.9473
$ws.Cells.Item(122, 9).Value = 4742.5  # I122: 4675 -> 4742.5
$ws.Cells.Item(122, 11).Value = 14227.5  # K122: 14025 -> 14227.5
$ws.Cells.Item(122, 13).Value = -11777.5  # M122: -11575 -> -11777.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3526.25  # H126: 3701.375 -> 3526.25
$ws.Cells.Item(126, 9).Value = 2490.75  # I126: 2802.3333 -> 2490.75
$ws.Cells.Item(126, 10).Value = 5597.25  # J126: 6398.5 -> 5597.25
$ws.Cells.Item(126, 11).Value = 7472.25  # K126: 8406.999899999999 -> 7472.25
$ws.Cells.Item(126, 12).Value = 16791.75  # L126: 19195.5 -> 16791.75
$ws.Cells.Item(126, 13).Value = -5002.25  # M126: -5936.999899999999 -> -5002.25
$ws.Cells.Item(126, 14).Value = -21731.75  # N126: -24135.5 -> -21731.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2833.7659  # H132: 2810.1667 -> 2833.7659
$ws.Cells.Item(132, 9).Value = 2923.2  # I132: 2812.0908 -> 2923.2
$ws.Cells.Item(132, 11).Value = 8769.599999999999  # K132: 8436.2724 -> 8769.599999999999
$ws.Cells.Item(132, 13).Value = -6239.599999999999  # M132: -5906.2724 -> -6239.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2532.5557  # H7: 2724.25 -> 2532.5557
$ws.Cells.Item(7, 9).Value = 1684.7142  # I7: 1799 -> 1684.7142
$ws.Cells.Item(7, 11).Value = 1684.7142  # K7: 1799 -> 1684.7142
$ws.Cells.Item(7, 13).Value = -1572.7142  # M7: -1687 -> -1572.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4137.7646  # H61: 4311.4707 -> 4137.7646
$ws.Cells.Item(61, 9).Value = 2326.3333  # I61: 2476.6667 -> 2326.3333
$ws.Cells.Item(61, 10).Value = 6175.625  # J61: 6375.625 -> 6175.625
$ws.Cells.Item(61, 11).Value = 2326.3333  # K61: 2476.6667 -> 2326.3333
$ws.Cells.Item(61, 12).Value = 6175.625  # L61: 6375.625 -> 6175.625
$ws.Cells.Item(61, 13).Value = -2124.3333  # M61: -2274.6667 -> -2124.3333
$ws.Cells.Item(61, 14).Value = -6579.625  # N61: -6779.625 -> -6579.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 4137.7646  # H113: 4311.4707 -> 4137.7646
$ws.Cells.Item(113, 9).Value = 2326.3333  # I113: 2476.6667 -> 2326.3333
$ws.Cells.Item(113, 10).Value = 6175.625  # J113: 6375.625 -> 6175.625
$ws.Cells.Item(113, 11).Value = 2326.3333  # K113: 2476.6667 -> 2326.3333
$ws.Cells.Item(113, 12).Value = 6175.625  # L113: 6375.625 -> 6175.625
$ws.Cells.Item(113, 13).Value = -156.3332999999998  # M113: -306.6667000000002 -> -156.3332999999998
$ws.Cells.Item(113, 14).Value = -10515.625  # N113: -10715.625 -> -10515.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 2532.5557  # H126: 2724.25 -> 2532.5557
$ws.Cells.Item(126, 9).Value = 1684.7142  # I126: 1799 -> 1684.7142
$ws.Cells.Item(126, 11).Value = 5054.142599999999  # K126: 5397 -> 5054.142599999999
$ws.Cells.Item(126, 13).Value = -2584.142599999999  # M126: -2927 -> -2584.142599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 0  # H49: 21000 -> 0
$ws.Cells.Item(49, 10).Value = 0  # J49: 21000 -> 0
$ws.Cells.Item(49, 12).Value = 0  # L49: 21000 -> 0
$ws.Cells.Item(49, 14).ClearContents()  # N49 remove (was -21460)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10750  # H62: 10565.218 -> 10750
$ws.Cells.Item(62, 10).Value = 12321.429  # J62: 11933.333 -> 12321.429
$ws.Cells.Item(62, 12).Value = 12321.429  # L62: 11933.333 -> 12321.429
$ws.Cells.Item(62, 14).Value = -13569.429  # N62: -13181.333 -> -13569.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 10750  # H65: 10565.218 -> 10750
$ws.Cells.Item(65, 10).Value = 12321.429  # J65: 11933.333 -> 12321.429
$ws.Cells.Item(65, 12).Value = 61607.145  # L65: 59666.665 -> 61607.145
$ws.Cells.Item(65, 14).Value = -67847.145  # N65: -65906.66500000001 -> -67847.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1268  # H113: 1224.6 -> 1268
$ws.Cells.Item(113, 9).Value = 1228.9  # I113: 1153.5454 -> 1228.9
$ws.Cells.Item(113, 11).Value = 3686.7  # K113: 3460.6362 -> 3686.7
$ws.Cells.Item(113, 13).Value = -1516.7  # M113: -1290.6362 -> -1516.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2657.5957  # H122: 2702.652 -> 2657.5957
$ws.Cells.Item(122, 9).Value = 1290  # I122: 1309.0541 -> 1290
$ws.Cells.Item(122, 11).Value = 3870  # K122: 3927.1623 -> 3870
$ws.Cells.Item(122, 13).Value = -1420  # M122: -1477.1623 -> -1420
